# Commit: "Documentation for language preparation done!"
#
# The author fixed a handful of translation keys/values in the
# "translation__am" sheet:
#   - three B-column values had accidentally picked up a leading `"`
#     character (email_address / opened_documents / sector_name rows);
#   - three A-column keys were renamed for clarity:
#       prepare_langauge_and_convert_to_JSON  -> prepare_language_and_convert_to_JSON
#       download_translation_template         -> download_translation_template_en
#       download_amharic_translation_template  -> download_translation_template_am

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the stray leading double-quote from the three Amharic values.
$ws.Range("B37").Value  = "የኢሜይል አድራሻ"
$ws.Range("B91").Value  = "ለአስተያየት ክፍት የሆኑ"
$ws.Range("B142").Value = "የሴክተሩ ስም"

# Rename the three translation-key identifiers (order matters for the
# resulting shared-string append order: en, then am, then prepare_language).
$ws.Range("A172").Value = "download_translation_template_en"
$ws.Range("A173").Value = "download_translation_template_am"
$ws.Range("A171").Value = "prepare_language_and_convert_to_JSON"

# Restore the cursor/selection position the author ended up on.
$ws.Range("A172").Select()
